# Generate Report for Handback
# Fills in the "Latest Handback File" / "Latest Handback DateTime" / "Error Detail"
# columns for the bcfc24a3-a958-4564-9432-cfeace43c908 row on both the zh-cn and
# de-de report sheets, and widens the "Error Detail" column.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1488dc921448870695948ea7719c37da7742b817/e2e/bcfc24a3-a958-4564-9432-cfeace43c908.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/620bce995be39ffb49db95d8f3fc8ec666eb0da2/e2e/bcfc24a3-a958-4564-9432-cfeace43c908.md."
$handbackUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/620bce995be39ffb49db95d8f3fc8ec666eb0da2/e2e/bcfc24a3-a958-4564-9432-cfeace43c908.md"
$handbackDisplay = "bcfc24a3-a958-4564-9432-cfeace43c908.md"

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("I8").Value = $handbackDisplay
$wsZh.Hyperlinks.Add($wsZh.Range("I8"), $handbackUrl, "", "", $handbackDisplay)
$wsZh.Range("J8").Value = "bcfc24a3-a958-4564-9432-cfeace43c908.ae48eef267a83ad1b505b75899aa391e89ee8c37.zh-cn.xlf"
$wsZh.Range("K8").Value = "2016-09-05 16:50:39"
$wsZh.Range("P8").Value = $errorDetail
$wsZh.Columns.Item(16).ColumnWidth = 39.17

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("I8").Value = $handbackDisplay
$wsDe.Hyperlinks.Add($wsDe.Range("I8"), $handbackUrl, "", "", $handbackDisplay)
$wsDe.Range("J8").Value = "bcfc24a3-a958-4564-9432-cfeace43c908.ae48eef267a83ad1b505b75899aa391e89ee8c37.de-de.xlf"
$wsDe.Range("K8").Value = "2016-09-05 16:50:47"
$wsDe.Range("P8").Value = $errorDetail
$wsDe.Columns.Item(16).ColumnWidth = 39.17
